$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update teaching style names - drop the "(offline)" / "(online)" suffixes.
# A leading apostrophe preserves the existing quotePrefix cell style instead
# of Excel silently re-deriving a fresh (slightly different) style.
$ws.Range("B2").Value = "'Truc tiep"
$ws.Range("B3").Value = "'Truc tuyen"

# Shrink column B to fit the new (shorter) text. The sheet was authored with
# a manually "best fit" column width (17.90625 -> 10.36328125 characters);
# reproduce that narrower width here.
$ws.Columns.Item(2).ColumnWidth = 9.5

# Move the active selection to C8
$ws.Range("C8").Select() | Out-Null
